$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 72, shifting existing rows 72:166 down to 73:167
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new record's data
$ws.Cells.Item(72, 1).Value = 3
$ws.Cells.Item(72, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(72, 3).Value = "Coquimbo"
$ws.Cells.Item(72, 4).Value = 44413
$ws.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(72, 5).Value = 5
$ws.Cells.Item(72, 6).Value = 100112043
$ws.Cells.Item(72, 7).Value = "Pepino ensalada"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 110
$ws.Cells.Item(72, 11).Value = 15000
$ws.Cells.Item(72, 12).Value = 16000
$ws.Cells.Item(72, 13).Value = 15545
$ws.Cells.Item(72, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 222
$ws.Cells.Item(72, 17).Value = 70
$ws.Cells.Item(72, 18).Value = "Hortaliza"
